# Apply the equity-history update to Sheet1:
#  1. Re-format column I ("trade_date") for the existing data rows (2-393) so it
#     uses the same datetime number format as column B ("datetime") instead of
#     the plain date format it had before.
#  2. Append four new trading-day rows (394-397) with the same column layout,
#     where the new "datetime" (B) cell gets the datetime format and the new
#     "trade_date" (I) cell keeps the plain date format (i.e. the format column
#     I used to have before step 1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$datetimeFormat = "YYYY-MM-DD HH:MM:SS"
$dateFormat = "YYYY-MM-DD"

$lastRow = $ws.UsedRange.Rows.Count

# 1. Re-stamp column I's number format for every existing data row.
$ws.Range("I2:I$lastRow").NumberFormat = $datetimeFormat

# 2. Append the new rows.
$newRows = @(
    @{ A = 14.84; B = 46049; C = "NSE"; D = 15;    E = 14.42; F = 14.8;  G = "JAIPOW"; H = 43748379; I = 46049; J = "INE351F01018"; K = "Jaiprakash Power Ventures Ltd"; L = "JAIPOW"; M = "BREEZE" },
    @{ A = 15.83; B = 46050; C = "NSE"; D = 15.94; E = 14.98; F = 14.98; G = "JAIPOW"; H = 52396279; I = 46050; J = "INE351F01018"; K = "Jaiprakash Power Ventures Ltd"; L = "JAIPOW"; M = "BREEZE" },
    @{ A = 14.88; B = 46051; C = "NSE"; D = 15.97; E = 14.85; F = 15.9;  G = "JAIPOW"; H = 49105845; I = 46051; J = "INE351F01018"; K = "Jaiprakash Power Ventures Ltd"; L = "JAIPOW"; M = "BREEZE" },
    @{ A = 15.19; B = 46052; C = "NSE"; D = 15.28; E = 14.49; F = 14.85; G = "JAIPOW"; H = 55817590; I = 46052; J = "INE351F01018"; K = "Jaiprakash Power Ventures Ltd"; L = "JAIPOW"; M = "BREEZE" }
)

$row = $lastRow
foreach ($data in $newRows) {
    $row = $row + 1

    $ws.Cells.Item($row, 1).Value = $data.A
    $ws.Cells.Item($row, 2).Value = $data.B
    $ws.Cells.Item($row, 2).NumberFormat = $datetimeFormat
    $ws.Cells.Item($row, 3).Value = $data.C
    $ws.Cells.Item($row, 4).Value = $data.D
    $ws.Cells.Item($row, 5).Value = $data.E
    $ws.Cells.Item($row, 6).Value = $data.F
    $ws.Cells.Item($row, 7).Value = $data.G
    $ws.Cells.Item($row, 8).Value = $data.H
    $ws.Cells.Item($row, 9).Value = $data.I
    $ws.Cells.Item($row, 9).NumberFormat = $dateFormat
    $ws.Cells.Item($row, 10).Value = $data.J
    $ws.Cells.Item($row, 11).Value = $data.K
    $ws.Cells.Item($row, 12).Value = $data.L
    $ws.Cells.Item($row, 13).Value = $data.M
}
